$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.000000000000000002092811143995738
$ws.Range("C3").Value = 0.8422882724052781
$ws.Range("C4").Value = 0.00003547663963132689
$ws.Range("C5").Value = 0.0000000000000000009776561004032234
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0.000000000000000003642974848100913
$ws.Range("C8").Value = 0.000000000000000003177071318118062
$ws.Range("C9").Value = 0.00219714938765768
$ws.Range("C10").Value = 0.006090650936498306
$ws.Range("C11").Value = 0.0003074966821064281
$ws.Range("C12").Value = 0.000000000000000000005283386679673052
$ws.Range("C13").Value = 0.149080953948828
